$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("reviews_count") is entirely empty of data (header only).
# Delete the entire column, shifting F:K (reviews_average..latest_review_date)
# left into E:J, and updating the used range/dimension accordingly.
$ws.Range("E:E").EntireColumn.Delete()
